$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.2594103048008066
$ws.Range("J2").Value = 0.2594103048008066
$ws.Range("M2").Value = 3.135398666666667
$ws.Range("N2").Value = 9.406196000000001
$ws.Range("O2").Value = 0.1723049126704688
$ws.Range("P2").Value = 0.1723049126704688
$ws.Range("Q2").Value = 0.04194849876133334
$ws.Range("R2").Value = 0.377536488852
$ws.Range("S2").Value = 0.04469766991452268
$ws.Range("T2").Value = 0.04469766991452268
$ws.Range("I3").Value = 0.2594103048008066
$ws.Range("J3").Value = 0.2594103048008066
$ws.Range("O3").Value = 0.1733096678828815
$ws.Range("P3").Value = 0.1733096678828815
$ws.Range("S3").Value = 0.04495831377042485
$ws.Range("T3").Value = 0.04495831377042483
$ws.Range("I4").Value = 0.2594103048008066
$ws.Range("J4").Value = 0.2594103048008066
$ws.Range("M4").Value = 0.4900660000000001
$ws.Range("N4").Value = 1.470198
$ws.Range("O4").Value = 0.02693143306797965
$ws.Range("P4").Value = 0.02693143306797965
$ws.Range("Q4").Value = 0.006556593014000001
$ws.Range("R4").Value = 0.059009337126
$ws.Range("S4").Value = 0.006986291260887124
$ws.Range("T4").Value = 0.006986291260887122
$ws.Range("I5").Value = 0.2594103048008066
$ws.Range("J5").Value = 0.2594103048008066
$ws.Range("M5").Value = 11.417657
$ws.Range("N5").Value = 34.252971
$ws.Range("O5").Value = 0.62745398637867
$ws.Range("P5").Value = 0.6274539863786701
$ws.Range("Q5").Value = 0.152756833003
$ws.Range("R5").Value = 1.374811497027
$ws.Range("S5").Value = 0.162768029854972
$ws.Range("T5").Value = 0.1627680298549719
$ws.Range("G6").Value = 0.03819566666666666
$ws.Range("H6").Value = 0.114587
$ws.Range("I6").Value = 0.7405896951991934
$ws.Range("J6").Value = 0.7405896951991934
$ws.Range("M6").Value = 3.135398666666667
$ws.Range("N6").Value = 9.406196000000001
$ws.Range("O6").Value = 0.1723049126704688
$ws.Range("P6").Value = 0.1723049126704688
$ws.Range("Q6").Value = 0.1197586423391111
$ws.Range("R6").Value = 1.077827781052
$ws.Range("S6").Value = 0.1276072427559461
$ws.Range("T6").Value = 0.1276072427559461
$ws.Range("G7").Value = 0.03819566666666666
$ws.Range("H7").Value = 0.114587
$ws.Range("I7").Value = 0.7405896951991934
$ws.Range("J7").Value = 0.7405896951991934
$ws.Range("O7").Value = 0.1733096678828815
$ws.Range("P7").Value = 0.1733096678828815
$ws.Range("Q7").Value = 0.1204569864446667
$ws.Range("R7").Value = 1.084112878002
$ws.Range("S7").Value = 0.1283513541124566
$ws.Range("T7").Value = 0.1283513541124566
$ws.Range("G8").Value = 0.03819566666666666
$ws.Range("H8").Value = 0.114587
$ws.Range("I8").Value = 0.7405896951991934
$ws.Range("J8").Value = 0.7405896951991934
$ws.Range("M8").Value = 0.4900660000000001
$ws.Range("N8").Value = 1.470198
$ws.Range("O8").Value = 0.02693143306797965
$ws.Range("P8").Value = 0.02693143306797965
$ws.Range("Q8").Value = 0.01871839758066667
$ws.Range("R8").Value = 0.168465578226
$ws.Range("S8").Value = 0.01994514180709253
$ws.Range("T8").Value = 0.01994514180709253
$ws.Range("G9").Value = 0.03819566666666666
$ws.Range("H9").Value = 0.114587
$ws.Range("I9").Value = 0.7405896951991934
$ws.Range("J9").Value = 0.7405896951991934
$ws.Range("M9").Value = 11.417657
$ws.Range("N9").Value = 34.252971
$ws.Range("O9").Value = 0.62745398637867
$ws.Range("P9").Value = 0.6274539863786701
$ws.Range("Q9").Value = 0.4361050208863333
$ws.Range("R9").Value = 3.924945187977
$ws.Range("S9").Value = 0.464685956523698
$ws.Range("T9").Value = 0.4646859565236982
